$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("priorities")

# Prepend "fisheries sensitive watershed" classification notes to the
# habitat comments for the rows that fall inside the newly classified
# FSW (rows 6, 19, 23, 29, 31 -> H6, H19, H23, H29, H31).

$ws.Range("H6").Value = "Classified as fisheries sensitive watershed ander FRPA for Bull Trout and Arctic Grayling (Beaudry 2013, FSW-TAG f-7-020). Some deep pools for overwintering and rearing.  Large woody debris and undercut banks throughout. Sections of gravel suitable for spawning.  Good flow.  "

$ws.Range("H29").Value = "Classified as fisheries sensitive watershed ander FRPA for Bull Trout and Arctic Grayling (Beaudry 2013, FSW-TAG f-7-020). Larger stream with good flow and high habitat complexity.  Frequent pockets of gravel suitable for spawing at pool tailouts and behind large woody debris. "

$ws.Range("H31").Value = "Classified as fisheries sensitive watershed ander FRPA for Bull Trout and Arctic Grayling (Beaudry 2013, FSW-TAG f-7-020). Large woody debris and pools throughout.  Frequent pockets of gravel suitable for spawning. "

$ws.Range("H19").Value = " Classified as fisheries sensitive watershed ander FRPA for Bull Trout and Arctic Grayling (Beaudry 2014, FSW-TAG f-7-022). Some deep pools and boulders, udercut banks,  large wody debris and gravels throughout.   Some debris steps from 30 - 70 cms high. Passble railway culvert located downstream (16603641). New bridge upstream."

$ws.Range("H23").Value = "Classified as fisheries sensitive watershed ander FRPA for Bull Trout and Arctic Grayling (Beaudry 2014, FSW-TAG f-7-022). CN Rail crossing.  Abundant gravels, large woody debris, undercut banks, overhanging vegetation and small woody debris. Recently installed bridges downsteam and upstream.  20 cm long bull trout (suspected) observed approximately 340 m upstream of the culvert.  Minnowtrapping conducted upstream and downstream with Rainbow Trout captured downstream."

# Leave the cursor where the author last left it before saving.
$ws.Range("D20").Select()
